# Invoice template cleanup: remove the pre-filled line-total formulas from
# the empty item rows (H16:H30) so each row starts blank instead of showing
# "=ItemQty*ItemUnitPrice" (which evaluates to 0 until the row is filled in),
# then leave the selection on the first line-item cell (H16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the formulas (but keep the existing cell formatting/style) from the
# 15 line-item rows. ClearContents() removes the formula/value while
# preserving number formatting, fill, borders, etc.
$ws.Range("H16:H30").ClearContents()

# Leave the cursor on the first item row's total cell.
$ws.Range("H16").Select()
